$wb = $excel.ActiveWorkbook

# Sheet "展览" (sheet1) - data sheet for exhibition listings
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value = 569
$ws1.Range("F4").Value = 1144
$ws1.Range("F5").Value = 127
$ws1.Range("F6").Value = 93
$ws1.Range("F9").Value = 1177
$ws1.Range("F10").Value = 16629
$ws1.Range("F11").Value = 291
$ws1.Range("F12").Value = 210
$ws1.Range("F13").Value = 1046
$ws1.Range("F14").Value = 6444
$ws1.Range("F16").Value = 133
$ws1.Range("F31").Value = 5074
$ws1.Range("F32").Value = 514
$ws1.Range("F33").Value = 11448
$ws1.Range("F36").Value = 164
$ws1.Range("F37").Value = 223
$ws1.Range("F38").Value = 3856
$ws1.Range("F39").Value = 273

# Sheet "全部类型" (sheet4) - combined listing with one extra row, so F indices are shifted by 1 for rows >= 33
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F3").Value = 569
$ws4.Range("F4").Value = 1144
$ws4.Range("F5").Value = 127
$ws4.Range("F6").Value = 93
$ws4.Range("F9").Value = 1177
$ws4.Range("F10").Value = 16629
$ws4.Range("F11").Value = 291
$ws4.Range("F12").Value = 210
$ws4.Range("F13").Value = 1046
$ws4.Range("F14").Value = 6444
$ws4.Range("F16").Value = 133
$ws4.Range("F31").Value = 5074
$ws4.Range("F32").Value = 514
$ws4.Range("F34").Value = 11448
$ws4.Range("F37").Value = 164
$ws4.Range("F38").Value = 223
$ws4.Range("F39").Value = 3856
$ws4.Range("F40").Value = 273
